# ArtifactRandomizeByLevel.xlsx - update level-based randomization table
# Commit message: "added player equip item"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artifact")

# Updated values for columns A (level), B (required stat), D (weight) per row.
# Column C and E are unchanged except C10 (0.9 -> 1).
$data = @(
    @{ Row = 2;  A = 5;  B = 20;  C = 0.1; D = 0.1; E = 10 },
    @{ Row = 3;  A = 10; B = 50;  C = 0.2; D = 0.2; E = 20 },
    @{ Row = 4;  A = 15; B = 90;  C = 0.3; D = 0.2; E = 30 },
    @{ Row = 5;  A = 20; B = 130; C = 0.4; D = 0.2; E = 40 },
    @{ Row = 6;  A = 25; B = 160; C = 0.5; D = 0.3; E = 50 },
    @{ Row = 7;  A = 35; B = 200; C = 0.6; D = 0.3; E = 60 },
    @{ Row = 8;  A = 45; B = 250; C = 0.7; D = 0.4; E = 70 },
    @{ Row = 9;  A = 55; B = 300; C = 0.8; D = 0.4; E = 80 },
    @{ Row = 10; A = 65; B = 400; C = 1;   D = 0.5; E = 90 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
}

# Update the active selection to match the saved workbook state (F6).
$ws.Range("F6").Select()
